# "details added to excel file"
# The repository link for the third row (C3) is filled in with the new
# project's repo URL, turned into a real hyperlink, and styled to match
# the other hyperlink cells (B2, C2, B3) which use the built-in
# "Hyperlink" cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRepoUrl = "https://github.com/dohaqabeel/OS_Project.git"

$target = $ws.Range("C3")
$target.Value = $newRepoUrl

# Create the actual OOXML hyperlink relationship for C3.
$null = $ws.Hyperlinks.Add($target, $newRepoUrl)

# Hyperlinks.Add() stamps its own ad-hoc font formatting onto the cell;
# re-apply the workbook's named "Hyperlink" style so C3 matches the
# existing hyperlink cells exactly (same style as B3/B2/C2).
$target.Style = "Hyperlink"

# Leave the cursor where the author left it when they saved the file.
$null = $ws.Range("B8").Select()
